# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Abril de 2020 a las 23:22"

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 844440
$ws.Range("C4").Value = 25696
$ws.Range("D4").Value = 83759
$ws.Range("E4").Value = 713454
$ws.Range("F4").Value = 14014
$ws.Range("G4").Value = 1909
$ws.Range("H4").Value = 47227

# --- India (row 20) ---
$ws.Range("B20").Value = 21370
$ws.Range("C20").Value = 1290
$ws.Range("D20").Value = 4370
$ws.Range("E20").Value = 16319
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 36
$ws.Range("H20").Value = 681

# --- Insert "Gabon" with fresh data above "Congo" / "Martinica" ---
# Row 129 previously held Congo's data, row 130 Martinica's data, row 131 Gabon's
# (stale) data. Gabon now gets newly updated figures and moves up to row 129,
# while Congo and Martinica keep their previous (unrefreshed) figures and shift
# down by one row each.
$ws.Range("A129").Value = "Gabon"
$ws.Range("B129").Value = 166
$ws.Range("C129").Value = 10
$ws.Range("D129").Value = 24
$ws.Range("E129").Value = 141
$ws.Range("F129").Value = 2
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 1

$ws.Range("A130").Value = "Congo"
$ws.Range("B130").Value = 165
$ws.Range("C130").Value = 0
$ws.Range("D130").Value = 16
$ws.Range("E130").Value = 143
$ws.Range("F130").Value = 0
$ws.Range("G130").Value = 0
$ws.Range("H130").Value = 6

$ws.Range("A131").Value = "Martinica"
$ws.Range("B131").Value = 164
$ws.Range("C131").Value = 1
$ws.Range("D131").Value = 73
$ws.Range("E131").Value = 77
$ws.Range("F131").Value = 6
$ws.Range("G131").Value = 0
$ws.Range("H131").Value = 14
